$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date updated
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was empty -> Alvearie Team
$ws.Range("B9").Value = "Alvearie Team"

# Row 10 ("Contact" / "No display for ContactDetail") -> Jurisdiction / United States of America
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row - remove it entirely,
# shifting all following rows up by one.
$ws.Rows.Item(11).Delete()

# "Case Sensitive" value was empty -> the literal text "true" (not a boolean).
# Plain Range.Value coerces "true"/"false" strings into real booleans, so stage the
# text in a scratch cell (using a leading apostrophe to force text), copy it, and
# paste-special (values only) into the target so it lands as a text/shared-string cell.
$helper = $ws.Range("B7")
$helper.Value = "'true"
$helper.Copy()
$ws.Range("B14").PasteSpecial(-4163)   # xlPasteValues
$helper.ClearContents()

# Restore B7's original formatting (the value-only paste above didn't touch it, but the
# apostrophe-prefixed text entry switched its number format to quoted text) by re-pasting
# the plain formatting from an untouched, identically-styled neighbor cell.
$fmtSource = $ws.Range("B9")
$fmtSource.Copy()
$helper.PasteSpecial(-4122)            # xlPasteFormats

$excel.CutCopyMode = 0

Write-Output "Metadata!B3  = $($ws.Range('B3').Value2)"
Write-Output "Metadata!B8  = $($ws.Range('B8').Value2)"
Write-Output "Metadata!B9  = $($ws.Range('B9').Value2)"
Write-Output "Metadata!A10 = $($ws.Range('A10').Value2)"
Write-Output "Metadata!B10 = $($ws.Range('B10').Value2)"
Write-Output "Metadata!A11 = $($ws.Range('A11').Value2)"
Write-Output "Metadata!B14 = $($ws.Range('B14').Value2)"
Write-Output "Metadata!B7  = $($ws.Range('B7').Value2)"
Write-Output "Metadata!dimension last row = $($ws.Range('A21').Value2)"
